$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - numbers 1..6
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 6

# Row 3 - letters a..f
$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "c"
$ws.Range("D3").Value = "d"
$ws.Range("E3").Value = "e"
$ws.Range("F3").Value = "f"

# Row 4 - numbers 7..12
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 12

# Row 5 - letters g..l
$ws.Range("A5").Value = "g"
$ws.Range("B5").Value = "h"
$ws.Range("C5").Value = "i"
$ws.Range("D5").Value = "j"
$ws.Range("E5").Value = "k"
$ws.Range("F5").Value = "l"

# Row 6 - numbers 13..18
$ws.Range("A6").Value = 13
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = 17
$ws.Range("F6").Value = 18

# Row 7 - letters m..r
$ws.Range("A7").Value = "m"
$ws.Range("B7").Value = "n"
$ws.Range("C7").Value = "o"
$ws.Range("D7").Value = "p"
$ws.Range("E7").Value = "q"
$ws.Range("F7").Value = "r"

# Row 8 - clear out the old values (A8=7, B8="g")
$ws.Range("A8:F8").ClearContents()

# Update the active cell selection to F7
$ws.Range("F7").Select()
